$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("time_variants")

# Clear the stray value in BC2 (keep formatting/style, just remove the value)
$ws.Range("BC2").ClearContents()

# Insert a new row above row 15; this shifts old rows 15-25 down to 16-26
$ws.Rows(15).Insert()

# Populate the newly inserted row 15 with the new parameter
$ws.Range("A15").Value = "program_perc_ipt_age0to5"
$ws.Range("B15").Value = "no"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = "yes"
$ws.Range("AW15").Value = 0
$ws.Range("BB15").Value = 0
$ws.Range("BC15").Value = 80
